# Update countries & provincias Spain
# - Re-sorted a handful of countries whose case counts changed enough to
#   move their rank in the (case-count-sorted) table: Senegal now ranks
#   above San Marino / Estado de Palestina; Mayotte now ranks above
#   Mauricio / Somalia; Madagascar now ranks above Camboya.
# - Refreshed the day's case/death/recovered figures for the affected
#   rows plus a handful of unrelated rows whose stats simply changed.
# - Bumped the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados..." timestamp (row 1) ---------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 13:22"

# --- Plain stat refreshes (no re-sort involved) ---------------------------

# Row 14: Brasil
$ws.Range("B14").Value = 50230
$ws.Range("C14").Value = 738
$ws.Range("E14").Value = 20314
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = 3343

# Row 37: Catar
$ws.Range("D37").Value = 809
$ws.Range("E37").Value = 7706

# Row 64: Kazajistan
$ws.Range("E64").Value = 1750
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 24

# Row 78: Eslovenia
$ws.Range("B78").Value = 1373
$ws.Range("C78").Value = 7
$ws.Range("E78").Value = 1082
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 80

# Row 123: Tanzania
$ws.Range("D123").Value = 48
$ws.Range("E123").Value = 226

# --- Senegal moves above San Marino / Estado de Palestina -----------------
# (rows 104-106 now show Senegal, San Marino, Estado de Palestina in turn,
#  each carrying its own refreshed stats)

$ws.Range("A104").Value = "Senegal"
$ws.Range("B104").Value = 545
$ws.Range("C104").Value = 66
$ws.Range("D104").Value = 262
$ws.Range("E104").Value = 277
$ws.Range("F104").Value = 1
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 6

$ws.Range("A105").Value = "San Marino"
$ws.Range("B105").Value = 501
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 63
$ws.Range("E105").Value = 398
$ws.Range("F105").Value = 3
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 40

$ws.Range("A106").Value = "Estado de Palestina"
$ws.Range("B106").Value = 480
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 92
$ws.Range("E106").Value = 384
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 4

# --- Mayotte moves above Mauricio / Somalia --------------------------------
# (rows 115-117 now show Mayotte, Mauricio, Somalia in turn)

$ws.Range("A115").Value = "Mayotte"
$ws.Range("B115").Value = 354
$ws.Range("C115").Value = 28
$ws.Range("D115").Value = 144
$ws.Range("E115").Value = 206
$ws.Range("F115").Value = 4
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 4

$ws.Range("A116").Value = "Mauricio"
$ws.Range("B116").Value = 331
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 266
$ws.Range("E116").Value = 56
$ws.Range("F116").Value = 3
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 9

$ws.Range("A117").Value = "Somalia"
$ws.Range("B117").Value = 328
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 8
$ws.Range("E117").Value = 304
$ws.Range("F117").Value = 2
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 16

# --- Madagascar moves above Camboya ----------------------------------------
# (rows 138-139 now show Madagascar, Camboya in turn)

$ws.Range("A138").Value = "Madagascar"
$ws.Range("B138").Value = 122
$ws.Range("C138").Value = 1
$ws.Range("D138").Value = 61
$ws.Range("E138").Value = 61
$ws.Range("F138").Value = 1
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0

$ws.Range("A139").Value = "Camboya"
$ws.Range("B139").Value = 122
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 110
$ws.Range("E139").Value = 12
$ws.Range("F139").Value = 1
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 0
